$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 10465
$ws1.Range("F3").Value = 434
$ws1.Range("F4").Value = 2533
$ws1.Range("F5").Value = 24
$ws1.Range("F6").Value = 286
$ws1.Range("F9").Value = 780
$ws1.Range("F12").Value = 1104
$ws1.Range("F13").Value = 3250
$ws1.Range("F14").Value = 2417
$ws1.Range("F15").Value = 61
$ws1.Range("F16").Value = 2171
$ws1.Range("F17").Value = 2171
$ws1.Range("F18").Value = 245
$ws1.Range("F19").Value = 1943
$ws1.Range("F22").Value = 581
$ws1.Range("F23").Value = 69
$ws1.Range("F24").Value = 257
$ws1.Range("F25").Value = 10
$ws1.Range("F26").Value = 27
$ws1.Range("F27").Value = 246
$ws1.Range("F28").Value = 51
$ws1.Range("F29").Value = 387
$ws1.Range("F30").Value = 10
$ws1.Range("F32").Value = 402
$ws1.Range("F33").Value = 606
$ws1.Range("F34").Value = 25
$ws1.Range("F35").Value = 56
$ws1.Range("F36").Value = 269
$ws1.Range("F37").Value = 10
$ws1.Range("F39").Value = 493
$ws1.Range("F40").Value = 473
$ws1.Range("F41").Value = 1725
$ws1.Range("F42").Value = 144
$ws1.Range("F43").Value = 448
$ws1.Range("F44").Value = 56
$ws1.Range("F45").Value = 469
$ws1.Range("F46").Value = 1048

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 29

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 10465
$ws4.Range("F3").Value = 434
$ws4.Range("F4").Value = 2533
$ws4.Range("F5").Value = 29
$ws4.Range("F6").Value = 24
$ws4.Range("F8").Value = 286
$ws4.Range("F11").Value = 780
$ws4.Range("F12").Value = 1104
$ws4.Range("F13").Value = 3250
$ws4.Range("F14").Value = 2417
$ws4.Range("F15").Value = 2171
$ws4.Range("F16").Value = 2171
$ws4.Range("F18").Value = 581
$ws4.Range("F19").Value = 69
$ws4.Range("F20").Value = 257
$ws4.Range("F21").Value = 10
$ws4.Range("F22").Value = 27
$ws4.Range("F23").Value = 246
$ws4.Range("F24").Value = 51
$ws4.Range("F25").Value = 387
$ws4.Range("F26").Value = 10
$ws4.Range("F28").Value = 402
$ws4.Range("F29").Value = 606
$ws4.Range("F30").Value = 25
$ws4.Range("F34").Value = 56
$ws4.Range("F35").Value = 269
$ws4.Range("F37").Value = 493
$ws4.Range("F39").Value = 473
$ws4.Range("F40").Value = 1725
$ws4.Range("F41").Value = 144
$ws4.Range("F45").Value = 448
$ws4.Range("F46").Value = 56
$ws4.Range("F47").Value = 469
$ws4.Range("F48").Value = 1048

